$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (A2:B2) is removed entirely; remaining rows shift up.
$ws.Rows(2).Delete()

# The old A1:B2 merge no longer applies to the new single-row layout.
$ws.Range("A1:B1").UnMerge()

# Drop the centered/merged formatting so the cells fall back to the
# default (unstyled) look.
$ws.Range("A1:C1").ClearFormats()

# New row 1 contents: A1=1, B1="Mohamed" (shared string), C1=23
$ws.Range("A1").Value = 1
$ws.Range("B1").Value = "Mohamed"
$ws.Range("C1").Value = 23
